# HIREINBOX_18_Month_Forecast.xlsx edits
# - Summary: Monthly Burn (Post-Launch) R 404,000 -> R 344,000
# - Team & Salaries: Full-Stack Developer salary R 60,000 -> R 80,000
# - Team & Salaries: remove founder salary rows (CEO Simon Rubin, Co-CEO Shay Sinbeti)
# - Expense Breakdown: Full-Stack Developer monthly cost 60000 -> 80000
# - Expense Breakdown: remove CEO / Co-CEO salary line items
# - Expense Breakdown: TOTAL 404000 -> 344000

$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B9").Value = "R 344,000"

# --- Team & Salaries sheet ---
$wsTeam = $wb.Worksheets.Item("Team & Salaries")
$wsTeam.Range("B5").Value = "R 80,000"
# Remove the founder salary rows entirely (CEO Simon Rubin / Co-CEO Shay Sinbeti)
$wsTeam.Range("A7:D8").ClearContents()

# --- Expense Breakdown sheet ---
$wsExpense = $wb.Worksheets.Item("Expense Breakdown")
$wsExpense.Range("C5").Value = 80000
# Remove the CEO / Co-CEO salary line items
$wsExpense.Range("A7:C8").ClearContents()
$wsExpense.Range("C19").Value = 344000
